$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.912.38'
$ws.Range("E2").Value = '  +0.51%  '

# Row 3
$ws.Range("D3").Value = '2.387.53'
$ws.Range("E3").Value = '  +6.41%  '

# Row 4
$ws.Range("E4").Value = '  -0.40%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.05'
$ws.Range("E5").Value = '  +10.75%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.10'
$ws.Range("E6").Value = '  -7.57%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.646'
$ws.Range("E7").Value = '  +2.78%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.651'
$ws.Range("E9").Value = '  +7.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.51'
$ws.Range("E10").Value = '  -5.53%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0936'
$ws.Range("E11").Value = '  +1.51%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.65'
$ws.Range("E12").Value = '  -3.91%  '

# Row 13
$ws.Range("E13").Value = '  -1.70%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.21'
$ws.Range("E14").Value = '  +13.39%  '

# Row 15
$ws.Range("E15").Value = '  +1.62%  '

# Row 16
$ws.Range("D16").Value = '2.747.04'
$ws.Range("E16").Value = '  +6.31%  '

# Row 17
$ws.Range("D17").Value = '2.397.40'
$ws.Range("E17").Value = '  +5.48%  '

# Row 18
$ws.Range("D18").Value = '43.016.65'
$ws.Range("E18").Value = '  +0.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.84'
$ws.Range("E19").Value = '  +9.18%  '

# Row 20
$ws.Range("E20").Value = '  +1.91%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.51'
$ws.Range("E21").Value = '  +2.57%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.63'
$ws.Range("E22").Value = '  +4.18%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '273.18'
$ws.Range("E23").Value = '  +9.10%  '

# Row 24
$ws.Range("E24").Value = '  -0.50%  '

# Row 25
$ws.Range("E25").Value = '  +7.36%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.71'
$ws.Range("E26").Value = '  +1.49%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.98'
$ws.Range("E28").Value = '  +4.15%  '

# Row 29
$ws.Range("E29").Value = '  -1.74%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.88'
$ws.Range("E30").Value = '  -0.76%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.33'
$ws.Range("E31").Value = '  -0.55%  '

# Row 32
$ws.Range("E32").Value = '  -0.30%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0931'
$ws.Range("E33").Value = '  +4.78%  '

# Row 34
$ws.Range("E34").Value = '  +3.56%  '

# Row 35
$ws.Range("E35").Value = '  +5.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.92'
$ws.Range("E36").Value = '  -3.43%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.15'
$ws.Range("E37").Value = '  -2.45%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0367'
$ws.Range("E38").Value = '  -2.88%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.106'
$ws.Range("E39").Value = '  +1.95%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.80'
$ws.Range("E40").Value = '  +16.31%  '

# Row 41
$ws.Range("E41").Value = '  +19.53%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.233'
$ws.Range("E42").Value = '  +0.92%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.66'
$ws.Range("E43").Value = '  -3.21%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '121.52'
$ws.Range("E44").Value = '  +15.07%  '

# Row 45
$ws.Range("E45").Value = '  +0.07%  '

# Row 46
$ws.Range("B46").Value = 'BitcoinSV'
$ws.Range("C46").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.77'
$ws.Range("E46").Value = '  +53.10%  '

# Row 47
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.35'
$ws.Range("E47").Value = '  -0.81%  '

# Row 48
$ws.Range("E48").Value = '  +9.04%  '

# Row 49
$ws.Range("E49").Value = '  +0.21%  '

# Row 50
$ws.Range("E50").Value = '  +0.65%  '

# Row 51
$ws.Range("D51").Value = '1.592.14'
$ws.Range("E51").Value = '  +10.41%  '
